$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I (I0) and J (IF), added to the right of the
# existing table (which ran through column H). Match the header row's
# existing style (bold + thin border + centered) by copying formats from H1.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data values for columns I and J, rows 2-47 (row, I-value, J-value).
$rows = @(
    @(2,7,7),
    @(3,4,5),
    @(4,9,9),
    @(5,6,7),
    @(6,10,10),
    @(7,11,11),
    @(8,3,4),
    @(9,8,8),
    @(10,9,9),
    @(11,9,9),
    @(12,7,7),
    @(13,8,9),
    @(14,7,7),
    @(15,6,6),
    @(16,5,5),
    @(17,7,7),
    @(18,8,8),
    @(19,8,9),
    @(20,3,4),
    @(21,8,8),
    @(22,8,8),
    @(23,10,10),
    @(24,6,6),
    @(25,7,7),
    @(26,7,7),
    @(27,7,7),
    @(28,7,7),
    @(29,6,6),
    @(30,7,7),
    @(31,3,4),
    @(32,4,4),
    @(33,6,7),
    @(34,6,6),
    @(35,9,9),
    @(36,7,8),
    @(37,6,6),
    @(38,5,5),
    @(39,7,7),
    @(40,8,9),
    @(41,6,6),
    @(42,5,5),
    @(43,7,7),
    @(44,9,9),
    @(45,5,5),
    @(46,3,3),
    @(47,3,3)

)

foreach ($r in $rows) {
    $row = $r[0]
    $ws.Cells.Item($row, 9).Value = $r[1]
    $ws.Cells.Item($row, 10).Value = $r[2]
}
